$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Technology_selection")

# Header row
$ws.Range("A1").Value = "Technology"
$ws.Range("B1").Value = "Technology Selected"

# Technology rows with the selection flag
$ws.Range("A2").Value = "Solar_PV"
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = "P_Gas"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "HLR_Biomass"
$ws.Range("B4").Value = 1

# Mirror the author's UI state: Technology_selection tab active,
# selection resting on the cell right below the last entry.
$ws.Activate()
$ws.Range("B5").Select()
